$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

$ws.Range("A40").Value = 42513
$ws.Range("B40").Value = 132
$ws.Range("C40").Value = 124
$ws.Range("D40").Value = 0
$ws.Range("E40").Value = 8
$ws.Range("F40").Value = 124
$ws.Range("G40").Value = 0.93939393939393945
$ws.Range("H40").Value = 45.650378787990618
$ws.Range("I40").Value = 35.516666660550982
$ws.Range("J40").Value = 143.45000000554137
